$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("StrategyDictionaries")

# Update values in Sheet1 rows 2 and 3
$ws1.Range("F2").Value = "3m"
$ws1.Range("K2").ClearContents()

$ws1.Range("A3").Value = 2
$ws1.Range("C3").Value = "ETHUSDT"
$ws1.Range("F3").Value = "3m"
$ws1.Range("K3").ClearContents()

# Delete row 8 entirely
$ws1.Rows.Item(8).ClearContents()

# Update selection on Sheet1
$ws1.Range("C3").Select()

# Update the HA_VWAP strategy settings JSON in StrategyDictionaries sheet
$ws3.Range("B5").Value = '{"EMA": 200, "DistVWAP_PCT": 0.05, "Nb_Signals": 2, "ExitOnEmaCross": false}'

$wb.Save()
